$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) attualmente_positivi (sheet1): append rows 289:291 (dates 44497/44501/
#    44502 with their "n. attualmente positivi" counts), copying the
#    formatting of the last existing row (288) down first so the new cells
#    pick up the same date/number styles & borders.
# ---------------------------------------------------------------------------
$ws1.Range("A288:B288").Copy()
$ws1.Range("A289:B291").PasteSpecial(-4122)

$ws1.Cells.Item(289,1).Value = 44497
$ws1.Cells.Item(289,2).Value = 12
$ws1.Cells.Item(290,1).Value = 44501
$ws1.Cells.Item(290,2).Value = 13
$ws1.Cells.Item(291,1).Value = 44502
$ws1.Cells.Item(291,2).Value = 13

# ---------------------------------------------------------------------------
# 2) totale_contagiati (sheet2): append rows 273:275 in the same fashion.
# ---------------------------------------------------------------------------
$ws2.Range("A272:M272").Copy()
$ws2.Range("A273:M275").PasteSpecial(-4122)

$ws2.Cells.Item(273,1).Value = 44497
$ws2.Cells.Item(273,2).Value = 1286
$ws2.Cells.Item(273,3).Value = 1
$ws2.Cells.Item(273,4).Value = 48
$ws2.Cells.Item(273,5).Value = 94
$ws2.Cells.Item(273,6).Value = 151
$ws2.Cells.Item(273,7).Value = 157
$ws2.Cells.Item(273,8).Value = 212
$ws2.Cells.Item(273,9).Value = 234
$ws2.Cells.Item(273,10).Value = 173
$ws2.Cells.Item(273,11).Value = 126
$ws2.Cells.Item(273,12).Value = 68
$ws2.Cells.Item(273,13).Value = 23

$ws2.Cells.Item(274,1).Value = 44501
$ws2.Cells.Item(274,2).Value = 1287
$ws2.Cells.Item(274,3).Value = 1
$ws2.Cells.Item(274,4).Value = 48
$ws2.Cells.Item(274,5).Value = 94
$ws2.Cells.Item(274,6).Value = 151
$ws2.Cells.Item(274,7).Value = 157
$ws2.Cells.Item(274,8).Value = 212
$ws2.Cells.Item(274,9).Value = 235
$ws2.Cells.Item(274,10).Value = 173
$ws2.Cells.Item(274,11).Value = 126
$ws2.Cells.Item(274,12).Value = 68
$ws2.Cells.Item(274,13).Value = 23

$ws2.Cells.Item(275,1).Value = 44502
$ws2.Cells.Item(275,2).Value = 1287
$ws2.Cells.Item(275,3).Value = 0
$ws2.Cells.Item(275,4).Value = 48
$ws2.Cells.Item(275,5).Value = 94
$ws2.Cells.Item(275,6).Value = 151
$ws2.Cells.Item(275,7).Value = 157
$ws2.Cells.Item(275,8).Value = 212
$ws2.Cells.Item(275,9).Value = 235
$ws2.Cells.Item(275,10).Value = 173
$ws2.Cells.Item(275,11).Value = 126
$ws2.Cells.Item(275,12).Value = 68
$ws2.Cells.Item(275,13).Value = 23

# ---------------------------------------------------------------------------
# 3) Extend the two charts' series ranges so they cover the newly added
#    rows (attualmente_positivi!A3:B291, totale_contagiati!A3:C275).
# ---------------------------------------------------------------------------
$chart1 = $ws1.ChartObjects(1).Chart
$ser1 = $chart1.SeriesCollection(1)
$ser1.Formula = "=SERIES(attualmente_positivi!`$B`$1,attualmente_positivi!`$A`$3:`$A`$291,attualmente_positivi!`$B`$3:`$B`$291,1)"

$chart2 = $ws2.ChartObjects(1).Chart
$ser2a = $chart2.SeriesCollection(1)
$ser2a.Formula = "=SERIES(totale_contagiati!`$B`$1,totale_contagiati!`$A`$3:`$A`$275,totale_contagiati!`$B`$3:`$B`$275,1)"
$ser2b = $chart2.SeriesCollection(2)
$ser2b.Formula = "=SERIES(totale_contagiati!`$C`$1,totale_contagiati!`$A`$3:`$A`$275,totale_contagiati!`$C`$3:`$C`$275,2)"

# ---------------------------------------------------------------------------
# 4) Selection / active-sheet bookkeeping to match the author's saved view:
#    attualmente_positivi keeps cell C277 selected (no longer the active
#    tab), while totale_contagiati becomes the active tab with D275:M275
#    selected.
# ---------------------------------------------------------------------------
$ws1.Range("C277").Select()
$ws2.Activate()
$ws2.Range("D275:M275").Select()
